$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C & E: labels (shared strings) and raw data values ---
$ws.Range("C1").Value = "1-above_h"
$ws.Range("E1").Value = "1-above_h"
$ws.Range("C2").Value = 0.092
$ws.Range("E2").Value = 0.1625
$ws.Range("C3").Value = "1-above_l"
$ws.Range("E3").Value = "1-above_l"
$ws.Range("C4").Value = 0.084
$ws.Range("E4").Value = 0.15625
$ws.Range("C5").Value = "1-below_h"
$ws.Range("E5").Value = "1-below_h"
$ws.Range("C6").Value = 0.11875
$ws.Range("E6").Value = 0.134444444444
$ws.Range("C7").Value = "1-below_l"
$ws.Range("E7").Value = "1-below_l"
$ws.Range("C8").Value = 0.0738461538462
$ws.Range("E8").Value = 0.156666666667
$ws.Range("C9").Value = "2-above_h"
$ws.Range("E9").Value = "2-above_h"
$ws.Range("C10").Value = 0.0877777777778
$ws.Range("E10").Value = 0.117777777778
$ws.Range("C11").Value = "2-above_l"
$ws.Range("E11").Value = "2-above_l"
$ws.Range("C12").Value = 0.087
$ws.Range("E12").Value = 0.116666666667
$ws.Range("C13").Value = "2-below_h"
$ws.Range("E13").Value = "2-below_h"
$ws.Range("C14").Value = 0.0918181818182
$ws.Range("E14").Value = 0.152
$ws.Range("C15").Value = "2-below_l"
$ws.Range("E15").Value = "2-below_l"
$ws.Range("C16").Value = 0.0975
$ws.Range("E16").Value = 0.23125
$ws.Range("C17").Value = "3-above_h"
$ws.Range("E17").Value = "3-above_h"
$ws.Range("C18").Value = 0.0877777777778
$ws.Range("E18").Value = 0.1325
$ws.Range("C19").Value = "3-above_l"
$ws.Range("E19").Value = "3-above_l"
$ws.Range("C20").Value = 0.092
$ws.Range("E20").Value = 0.142307692308
$ws.Range("C21").Value = "3-below_h"
$ws.Range("E21").Value = "3-below_h"
$ws.Range("C22").Value = 0.154545454545
$ws.Range("E22").Value = 0.128888888889
$ws.Range("C23").Value = "3-below_l"
$ws.Range("E23").Value = "3-below_l"
$ws.Range("C24").Value = 0.15625
$ws.Range("E24").Value = 0.187692307692
$ws.Range("C25").Value = "5-above_h"
$ws.Range("E25").Value = "5-above_h"
$ws.Range("C26").Value = 0.0855555555556
$ws.Range("E26").Value = 0.087
$ws.Range("C27").Value = "5-above_l"
$ws.Range("E27").Value = "5-above_l"
$ws.Range("C28").Value = 0.0808333333333
$ws.Range("E28").Value = 0.130909090909
$ws.Range("C29").Value = "5-below_h"
$ws.Range("E29").Value = "5-below_h"
$ws.Range("C30").Value = 0.098
$ws.Range("E30").Value = 0.141
$ws.Range("C31").Value = "5-below_l"
$ws.Range("E31").Value = "5-below_l"
$ws.Range("C32").Value = 0.0915384615385
$ws.Range("E32").Value = 0.122
$ws.Range("C33").Value = "8-above_h"
$ws.Range("E33").Value = "8-above_h"
$ws.Range("C34").Value = 0.0875
$ws.Range("E34").Value = 0.086
$ws.Range("C35").Value = "8-above_l"
$ws.Range("E35").Value = "8-above_l"
$ws.Range("C36").Value = 0.081
$ws.Range("E36").Value = 0.0908333333333
$ws.Range("C37").Value = "8-below_h"
$ws.Range("E37").Value = "8-below_h"
$ws.Range("C38").Value = 0.0833333333333
$ws.Range("E38").Value = 0.13125
$ws.Range("C39").Value = "8-below_l"
$ws.Range("E39").Value = "8-below_l"
$ws.Range("C40").Value = 0.0972727272727
$ws.Range("E40").Value = 0.088

# --- Column D & F: AVERAGE formulas on summary rows ---
$ws.Range("D1").Formula = "=AVERAGE(C2,C4,C6,C8)"
$ws.Range("F1").Formula = "=AVERAGE(E2,E4,E6,E8)"
$ws.Range("D9").Formula = "=AVERAGE(C10,C12,C14,C16)"
$ws.Range("F9").Formula = "=AVERAGE(E10,E12,E14,E16)"
$ws.Range("D17").Formula = "=AVERAGE(C18,C20,C22,C24)"
$ws.Range("F17").Formula = "=AVERAGE(E18,E20,E22,E24)"
$ws.Range("D25").Formula = "=AVERAGE(C26,C28,C30,C32)"
$ws.Range("F25").Formula = "=AVERAGE(E26,E28,E30,E32)"
$ws.Range("D33").Formula = "=AVERAGE(C34,C36,C38,C40)"
$ws.Range("F33").Formula = "=AVERAGE(E34,E36,E38,E40)"

# --- Column D & F: number format 0.0000 across rows 1-40 (stamps style on blank cells too) ---
$ws.Range("D1:D40").NumberFormat = "0.0000"
$ws.Range("F1:F40").NumberFormat = "0.0000"

# --- Selection matches final state in diff ---
$ws.Range("F13").Select()
